$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "2023-06-28 Wednesday"

# Update each math-problem cell by exact table coordinates (row, column),
# in document order, to avoid any cross-cell text collisions.
$t = $d.Tables.Item(1)
$failed = @()

$c = $t.Cell(1,1)
if ($c.Range.Text.StartsWith("3+28=")) {
  $c.Range.Text = "99-83="
} else {
  $failed += "row 1 col 1: expected [3+28=] got [$($c.Range.Text)]"
}
$c = $t.Cell(1,2)
if ($c.Range.Text.StartsWith("76-36=")) {
  $c.Range.Text = "26+28="
} else {
  $failed += "row 1 col 2: expected [76-36=] got [$($c.Range.Text)]"
}
$c = $t.Cell(1,3)
if ($c.Range.Text.StartsWith("13+27=")) {
  $c.Range.Text = "43+38="
} else {
  $failed += "row 1 col 3: expected [13+27=] got [$($c.Range.Text)]"
}
$c = $t.Cell(1,4)
if ($c.Range.Text.StartsWith("50+47=")) {
  $c.Range.Text = "41+41="
} else {
  $failed += "row 1 col 4: expected [50+47=] got [$($c.Range.Text)]"
}
$c = $t.Cell(1,5)
if ($c.Range.Text.StartsWith("2+53=")) {
  $c.Range.Text = "1-0="
} else {
  $failed += "row 1 col 5: expected [2+53=] got [$($c.Range.Text)]"
}
$c = $t.Cell(2,1)
if ($c.Range.Text.StartsWith("49-39=")) {
  $c.Range.Text = "96-56="
} else {
  $failed += "row 2 col 1: expected [49-39=] got [$($c.Range.Text)]"
}
$c = $t.Cell(2,2)
if ($c.Range.Text.StartsWith("5+85=")) {
  $c.Range.Text = "18+2="
} else {
  $failed += "row 2 col 2: expected [5+85=] got [$($c.Range.Text)]"
}
$c = $t.Cell(2,3)
if ($c.Range.Text.StartsWith("35-31=")) {
  $c.Range.Text = "26+28="
} else {
  $failed += "row 2 col 3: expected [35-31=] got [$($c.Range.Text)]"
}
$c = $t.Cell(2,4)
if ($c.Range.Text.StartsWith("50-48=")) {
  $c.Range.Text = "47+19="
} else {
  $failed += "row 2 col 4: expected [50-48=] got [$($c.Range.Text)]"
}
$c = $t.Cell(2,5)
if ($c.Range.Text.StartsWith("34+60=")) {
  $c.Range.Text = "56-41="
} else {
  $failed += "row 2 col 5: expected [34+60=] got [$($c.Range.Text)]"
}
$c = $t.Cell(3,1)
if ($c.Range.Text.StartsWith("19-14=")) {
  $c.Range.Text = "96-69="
} else {
  $failed += "row 3 col 1: expected [19-14=] got [$($c.Range.Text)]"
}
$c = $t.Cell(3,2)
if ($c.Range.Text.StartsWith("22-18=")) {
  $c.Range.Text = "69-2="
} else {
  $failed += "row 3 col 2: expected [22-18=] got [$($c.Range.Text)]"
}
$c = $t.Cell(3,3)
if ($c.Range.Text.StartsWith("36-7=")) {
  $c.Range.Text = "96-83="
} else {
  $failed += "row 3 col 3: expected [36-7=] got [$($c.Range.Text)]"
}
$c = $t.Cell(3,4)
if ($c.Range.Text.StartsWith("29-13=")) {
  $c.Range.Text = "85-29="
} else {
  $failed += "row 3 col 4: expected [29-13=] got [$($c.Range.Text)]"
}
$c = $t.Cell(3,5)
if ($c.Range.Text.StartsWith("70-17=")) {
  $c.Range.Text = "73-56="
} else {
  $failed += "row 3 col 5: expected [70-17=] got [$($c.Range.Text)]"
}
$c = $t.Cell(4,1)
if ($c.Range.Text.StartsWith("65-3=")) {
  $c.Range.Text = "5+42="
} else {
  $failed += "row 4 col 1: expected [65-3=] got [$($c.Range.Text)]"
}
$c = $t.Cell(4,2)
if ($c.Range.Text.StartsWith("26-12=")) {
  $c.Range.Text = "55+10="
} else {
  $failed += "row 4 col 2: expected [26-12=] got [$($c.Range.Text)]"
}
$c = $t.Cell(4,3)
if ($c.Range.Text.StartsWith("34+4=")) {
  $c.Range.Text = "33-23="
} else {
  $failed += "row 4 col 3: expected [34+4=] got [$($c.Range.Text)]"
}
$c = $t.Cell(4,4)
if ($c.Range.Text.StartsWith("10+49=")) {
  $c.Range.Text = "46-8="
} else {
  $failed += "row 4 col 4: expected [10+49=] got [$($c.Range.Text)]"
}
$c = $t.Cell(4,5)
if ($c.Range.Text.StartsWith("56-0=")) {
  $c.Range.Text = "74-38="
} else {
  $failed += "row 4 col 5: expected [56-0=] got [$($c.Range.Text)]"
}
$c = $t.Cell(5,1)
if ($c.Range.Text.StartsWith("64+28=")) {
  $c.Range.Text = "11-0="
} else {
  $failed += "row 5 col 1: expected [64+28=] got [$($c.Range.Text)]"
}
$c = $t.Cell(5,2)
if ($c.Range.Text.StartsWith("30+53=")) {
  $c.Range.Text = "80-78="
} else {
  $failed += "row 5 col 2: expected [30+53=] got [$($c.Range.Text)]"
}
$c = $t.Cell(5,3)
if ($c.Range.Text.StartsWith("79-70=")) {
  $c.Range.Text = "8+11="
} else {
  $failed += "row 5 col 3: expected [79-70=] got [$($c.Range.Text)]"
}
$c = $t.Cell(5,4)
if ($c.Range.Text.StartsWith("78-68=")) {
  $c.Range.Text = "49+35="
} else {
  $failed += "row 5 col 4: expected [78-68=] got [$($c.Range.Text)]"
}
$c = $t.Cell(5,5)
if ($c.Range.Text.StartsWith("63+7=")) {
  $c.Range.Text = "19-14="
} else {
  $failed += "row 5 col 5: expected [63+7=] got [$($c.Range.Text)]"
}
$c = $t.Cell(6,1)
if ($c.Range.Text.StartsWith("42-18=")) {
  $c.Range.Text = "63-44="
} else {
  $failed += "row 6 col 1: expected [42-18=] got [$($c.Range.Text)]"
}
$c = $t.Cell(6,2)
if ($c.Range.Text.StartsWith("29+36=")) {
  $c.Range.Text = "81-55="
} else {
  $failed += "row 6 col 2: expected [29+36=] got [$($c.Range.Text)]"
}
$c = $t.Cell(6,3)
if ($c.Range.Text.StartsWith("39-17=")) {
  $c.Range.Text = "20+22="
} else {
  $failed += "row 6 col 3: expected [39-17=] got [$($c.Range.Text)]"
}
$c = $t.Cell(6,4)
if ($c.Range.Text.StartsWith("37+46=")) {
  $c.Range.Text = "75+3="
} else {
  $failed += "row 6 col 4: expected [37+46=] got [$($c.Range.Text)]"
}
$c = $t.Cell(6,5)
if ($c.Range.Text.StartsWith("40+47=")) {
  $c.Range.Text = "7+48="
} else {
  $failed += "row 6 col 5: expected [40+47=] got [$($c.Range.Text)]"
}
$c = $t.Cell(7,1)
if ($c.Range.Text.StartsWith("71-49=")) {
  $c.Range.Text = "51-40="
} else {
  $failed += "row 7 col 1: expected [71-49=] got [$($c.Range.Text)]"
}
$c = $t.Cell(7,2)
if ($c.Range.Text.StartsWith("0+84=")) {
  $c.Range.Text = "35+55="
} else {
  $failed += "row 7 col 2: expected [0+84=] got [$($c.Range.Text)]"
}
$c = $t.Cell(7,3)
if ($c.Range.Text.StartsWith("87-77=")) {
  $c.Range.Text = "1+3="
} else {
  $failed += "row 7 col 3: expected [87-77=] got [$($c.Range.Text)]"
}
$c = $t.Cell(7,4)
if ($c.Range.Text.StartsWith("19+79=")) {
  $c.Range.Text = "29+41="
} else {
  $failed += "row 7 col 4: expected [19+79=] got [$($c.Range.Text)]"
}
$c = $t.Cell(7,5)
if ($c.Range.Text.StartsWith("56+21=")) {
  $c.Range.Text = "26+50="
} else {
  $failed += "row 7 col 5: expected [56+21=] got [$($c.Range.Text)]"
}
$c = $t.Cell(8,1)
if ($c.Range.Text.StartsWith("0+65=")) {
  $c.Range.Text = "67-49="
} else {
  $failed += "row 8 col 1: expected [0+65=] got [$($c.Range.Text)]"
}
$c = $t.Cell(8,2)
if ($c.Range.Text.StartsWith("47-18=")) {
  $c.Range.Text = "9+28="
} else {
  $failed += "row 8 col 2: expected [47-18=] got [$($c.Range.Text)]"
}
$c = $t.Cell(8,3)
if ($c.Range.Text.StartsWith("70-66=")) {
  $c.Range.Text = "97-65="
} else {
  $failed += "row 8 col 3: expected [70-66=] got [$($c.Range.Text)]"
}
$c = $t.Cell(8,4)
if ($c.Range.Text.StartsWith("21+16=")) {
  $c.Range.Text = "16+72="
} else {
  $failed += "row 8 col 4: expected [21+16=] got [$($c.Range.Text)]"
}
$c = $t.Cell(8,5)
if ($c.Range.Text.StartsWith("49+20=")) {
  $c.Range.Text = "37+17="
} else {
  $failed += "row 8 col 5: expected [49+20=] got [$($c.Range.Text)]"
}
$c = $t.Cell(9,1)
if ($c.Range.Text.StartsWith("81-71=")) {
  $c.Range.Text = "24+16="
} else {
  $failed += "row 9 col 1: expected [81-71=] got [$($c.Range.Text)]"
}
$c = $t.Cell(9,2)
if ($c.Range.Text.StartsWith("67-5=")) {
  $c.Range.Text = "81+6="
} else {
  $failed += "row 9 col 2: expected [67-5=] got [$($c.Range.Text)]"
}
$c = $t.Cell(9,3)
if ($c.Range.Text.StartsWith("39-8=")) {
  $c.Range.Text = "9+73="
} else {
  $failed += "row 9 col 3: expected [39-8=] got [$($c.Range.Text)]"
}
$c = $t.Cell(9,4)
if ($c.Range.Text.StartsWith("78-49=")) {
  $c.Range.Text = "25+72="
} else {
  $failed += "row 9 col 4: expected [78-49=] got [$($c.Range.Text)]"
}
$c = $t.Cell(9,5)
if ($c.Range.Text.StartsWith("12+47=")) {
  $c.Range.Text = "1+63="
} else {
  $failed += "row 9 col 5: expected [12+47=] got [$($c.Range.Text)]"
}
$c = $t.Cell(10,1)
if ($c.Range.Text.StartsWith("0+27=")) {
  $c.Range.Text = "62+11="
} else {
  $failed += "row 10 col 1: expected [0+27=] got [$($c.Range.Text)]"
}
$c = $t.Cell(10,2)
if ($c.Range.Text.StartsWith("67-32=")) {
  $c.Range.Text = "36-6="
} else {
  $failed += "row 10 col 2: expected [67-32=] got [$($c.Range.Text)]"
}
$c = $t.Cell(10,3)
if ($c.Range.Text.StartsWith("17+50=")) {
  $c.Range.Text = "94+1="
} else {
  $failed += "row 10 col 3: expected [17+50=] got [$($c.Range.Text)]"
}
$c = $t.Cell(10,4)
if ($c.Range.Text.StartsWith("45+20=")) {
  $c.Range.Text = "85-8="
} else {
  $failed += "row 10 col 4: expected [45+20=] got [$($c.Range.Text)]"
}
$c = $t.Cell(10,5)
if ($c.Range.Text.StartsWith("47+24=")) {
  $c.Range.Text = "9+54="
} else {
  $failed += "row 10 col 5: expected [47+24=] got [$($c.Range.Text)]"
}
$c = $t.Cell(11,1)
if ($c.Range.Text.StartsWith("74+7=")) {
  $c.Range.Text = "67-26="
} else {
  $failed += "row 11 col 1: expected [74+7=] got [$($c.Range.Text)]"
}
$c = $t.Cell(11,2)
if ($c.Range.Text.StartsWith("38-19=")) {
  $c.Range.Text = "70-42="
} else {
  $failed += "row 11 col 2: expected [38-19=] got [$($c.Range.Text)]"
}
$c = $t.Cell(11,3)
if ($c.Range.Text.StartsWith("10+34=")) {
  $c.Range.Text = "27+26="
} else {
  $failed += "row 11 col 3: expected [10+34=] got [$($c.Range.Text)]"
}
$c = $t.Cell(11,4)
if ($c.Range.Text.StartsWith("49+21=")) {
  $c.Range.Text = "45+38="
} else {
  $failed += "row 11 col 4: expected [49+21=] got [$($c.Range.Text)]"
}
$c = $t.Cell(11,5)
if ($c.Range.Text.StartsWith("77+10=")) {
  $c.Range.Text = "84+5="
} else {
  $failed += "row 11 col 5: expected [77+10=] got [$($c.Range.Text)]"
}
$c = $t.Cell(12,1)
if ($c.Range.Text.StartsWith("8-3=")) {
  $c.Range.Text = "92-16="
} else {
  $failed += "row 12 col 1: expected [8-3=] got [$($c.Range.Text)]"
}
$c = $t.Cell(12,2)
if ($c.Range.Text.StartsWith("31+12=")) {
  $c.Range.Text = "79-2="
} else {
  $failed += "row 12 col 2: expected [31+12=] got [$($c.Range.Text)]"
}
$c = $t.Cell(12,3)
if ($c.Range.Text.StartsWith("28-13=")) {
  $c.Range.Text = "86-1="
} else {
  $failed += "row 12 col 3: expected [28-13=] got [$($c.Range.Text)]"
}
$c = $t.Cell(12,4)
if ($c.Range.Text.StartsWith("64-59=")) {
  $c.Range.Text = "14+75="
} else {
  $failed += "row 12 col 4: expected [64-59=] got [$($c.Range.Text)]"
}
$c = $t.Cell(12,5)
if ($c.Range.Text.StartsWith("2+73=")) {
  $c.Range.Text = "52-9="
} else {
  $failed += "row 12 col 5: expected [2+73=] got [$($c.Range.Text)]"
}
$c = $t.Cell(13,1)
if ($c.Range.Text.StartsWith("4+40=")) {
  $c.Range.Text = "28-5="
} else {
  $failed += "row 13 col 1: expected [4+40=] got [$($c.Range.Text)]"
}
$c = $t.Cell(13,2)
if ($c.Range.Text.StartsWith("58+20=")) {
  $c.Range.Text = "15+22="
} else {
  $failed += "row 13 col 2: expected [58+20=] got [$($c.Range.Text)]"
}
$c = $t.Cell(13,3)
if ($c.Range.Text.StartsWith("34+6=")) {
  $c.Range.Text = "41-7="
} else {
  $failed += "row 13 col 3: expected [34+6=] got [$($c.Range.Text)]"
}
$c = $t.Cell(13,4)
if ($c.Range.Text.StartsWith("56-23=")) {
  $c.Range.Text = "45+39="
} else {
  $failed += "row 13 col 4: expected [56-23=] got [$($c.Range.Text)]"
}
$c = $t.Cell(13,5)
if ($c.Range.Text.StartsWith("96-43=")) {
  $c.Range.Text = "94+1="
} else {
  $failed += "row 13 col 5: expected [96-43=] got [$($c.Range.Text)]"
}
$c = $t.Cell(14,1)
if ($c.Range.Text.StartsWith("45-36=")) {
  $c.Range.Text = "69+15="
} else {
  $failed += "row 14 col 1: expected [45-36=] got [$($c.Range.Text)]"
}
$c = $t.Cell(14,2)
if ($c.Range.Text.StartsWith("4+91=")) {
  $c.Range.Text = "34-11="
} else {
  $failed += "row 14 col 2: expected [4+91=] got [$($c.Range.Text)]"
}
$c = $t.Cell(14,3)
if ($c.Range.Text.StartsWith("91-67=")) {
  $c.Range.Text = "32+59="
} else {
  $failed += "row 14 col 3: expected [91-67=] got [$($c.Range.Text)]"
}
$c = $t.Cell(14,4)
if ($c.Range.Text.StartsWith("76-26=")) {
  $c.Range.Text = "63-15="
} else {
  $failed += "row 14 col 4: expected [76-26=] got [$($c.Range.Text)]"
}
$c = $t.Cell(14,5)
if ($c.Range.Text.StartsWith("1+41=")) {
  $c.Range.Text = "67-27="
} else {
  $failed += "row 14 col 5: expected [1+41=] got [$($c.Range.Text)]"
}
$c = $t.Cell(15,1)
if ($c.Range.Text.StartsWith("72-67=")) {
  $c.Range.Text = "85-2="
} else {
  $failed += "row 15 col 1: expected [72-67=] got [$($c.Range.Text)]"
}
$c = $t.Cell(15,2)
if ($c.Range.Text.StartsWith("93-72=")) {
  $c.Range.Text = "80+6="
} else {
  $failed += "row 15 col 2: expected [93-72=] got [$($c.Range.Text)]"
}
$c = $t.Cell(15,3)
if ($c.Range.Text.StartsWith("92-50=")) {
  $c.Range.Text = "98-94="
} else {
  $failed += "row 15 col 3: expected [92-50=] got [$($c.Range.Text)]"
}
$c = $t.Cell(15,4)
if ($c.Range.Text.StartsWith("18+59=")) {
  $c.Range.Text = "99-1="
} else {
  $failed += "row 15 col 4: expected [18+59=] got [$($c.Range.Text)]"
}
$c = $t.Cell(15,5)
if ($c.Range.Text.StartsWith("79-0=")) {
  $c.Range.Text = "21+73="
} else {
  $failed += "row 15 col 5: expected [79-0=] got [$($c.Range.Text)]"
}
$c = $t.Cell(16,1)
if ($c.Range.Text.StartsWith("1+49=")) {
  $c.Range.Text = "36+38="
} else {
  $failed += "row 16 col 1: expected [1+49=] got [$($c.Range.Text)]"
}
$c = $t.Cell(16,2)
if ($c.Range.Text.StartsWith("99-3=")) {
  $c.Range.Text = "13-6="
} else {
  $failed += "row 16 col 2: expected [99-3=] got [$($c.Range.Text)]"
}
$c = $t.Cell(16,3)
if ($c.Range.Text.StartsWith("58+35=")) {
  $c.Range.Text = "16+44="
} else {
  $failed += "row 16 col 3: expected [58+35=] got [$($c.Range.Text)]"
}
$c = $t.Cell(16,4)
if ($c.Range.Text.StartsWith("70-51=")) {
  $c.Range.Text = "27+46="
} else {
  $failed += "row 16 col 4: expected [70-51=] got [$($c.Range.Text)]"
}
$c = $t.Cell(16,5)
if ($c.Range.Text.StartsWith("45-6=")) {
  $c.Range.Text = "51+46="
} else {
  $failed += "row 16 col 5: expected [45-6=] got [$($c.Range.Text)]"
}
$c = $t.Cell(17,1)
if ($c.Range.Text.StartsWith("55+22=")) {
  $c.Range.Text = "10-4="
} else {
  $failed += "row 17 col 1: expected [55+22=] got [$($c.Range.Text)]"
}
$c = $t.Cell(17,2)
if ($c.Range.Text.StartsWith("52+44=")) {
  $c.Range.Text = "5+61="
} else {
  $failed += "row 17 col 2: expected [52+44=] got [$($c.Range.Text)]"
}
$c = $t.Cell(17,3)
if ($c.Range.Text.StartsWith("41+7=")) {
  $c.Range.Text = "57+21="
} else {
  $failed += "row 17 col 3: expected [41+7=] got [$($c.Range.Text)]"
}
$c = $t.Cell(17,4)
if ($c.Range.Text.StartsWith("95-90=")) {
  $c.Range.Text = "80+7="
} else {
  $failed += "row 17 col 4: expected [95-90=] got [$($c.Range.Text)]"
}
$c = $t.Cell(17,5)
if ($c.Range.Text.StartsWith("36-32=")) {
  $c.Range.Text = "52-10="
} else {
  $failed += "row 17 col 5: expected [36-32=] got [$($c.Range.Text)]"
}
$c = $t.Cell(18,1)
if ($c.Range.Text.StartsWith("27-26=")) {
  $c.Range.Text = "68-40="
} else {
  $failed += "row 18 col 1: expected [27-26=] got [$($c.Range.Text)]"
}
$c = $t.Cell(18,2)
if ($c.Range.Text.StartsWith("89-17=")) {
  $c.Range.Text = "39+28="
} else {
  $failed += "row 18 col 2: expected [89-17=] got [$($c.Range.Text)]"
}
$c = $t.Cell(18,3)
if ($c.Range.Text.StartsWith("76-65=")) {
  $c.Range.Text = "38+48="
} else {
  $failed += "row 18 col 3: expected [76-65=] got [$($c.Range.Text)]"
}
$c = $t.Cell(18,4)
if ($c.Range.Text.StartsWith("1+36=")) {
  $c.Range.Text = "30+5="
} else {
  $failed += "row 18 col 4: expected [1+36=] got [$($c.Range.Text)]"
}
$c = $t.Cell(18,5)
if ($c.Range.Text.StartsWith("91-58=")) {
  $c.Range.Text = "58+34="
} else {
  $failed += "row 18 col 5: expected [91-58=] got [$($c.Range.Text)]"
}
$c = $t.Cell(19,1)
if ($c.Range.Text.StartsWith("20+42=")) {
  $c.Range.Text = "6+19="
} else {
  $failed += "row 19 col 1: expected [20+42=] got [$($c.Range.Text)]"
}
$c = $t.Cell(19,2)
if ($c.Range.Text.StartsWith("14+46=")) {
  $c.Range.Text = "41+2="
} else {
  $failed += "row 19 col 2: expected [14+46=] got [$($c.Range.Text)]"
}
$c = $t.Cell(19,3)
if ($c.Range.Text.StartsWith("40-38=")) {
  $c.Range.Text = "75-15="
} else {
  $failed += "row 19 col 3: expected [40-38=] got [$($c.Range.Text)]"
}
$c = $t.Cell(19,4)
if ($c.Range.Text.StartsWith("71-11=")) {
  $c.Range.Text = "56+5="
} else {
  $failed += "row 19 col 4: expected [71-11=] got [$($c.Range.Text)]"
}
$c = $t.Cell(19,5)
if ($c.Range.Text.StartsWith("64-37=")) {
  $c.Range.Text = "95-20="
} else {
  $failed += "row 19 col 5: expected [64-37=] got [$($c.Range.Text)]"
}
$c = $t.Cell(20,1)
if ($c.Range.Text.StartsWith("1+20=")) {
  $c.Range.Text = "71-11="
} else {
  $failed += "row 20 col 1: expected [1+20=] got [$($c.Range.Text)]"
}
$c = $t.Cell(20,2)
if ($c.Range.Text.StartsWith("76-17=")) {
  $c.Range.Text = "37+45="
} else {
  $failed += "row 20 col 2: expected [76-17=] got [$($c.Range.Text)]"
}
$c = $t.Cell(20,3)
if ($c.Range.Text.StartsWith("88-17=")) {
  $c.Range.Text = "60-48="
} else {
  $failed += "row 20 col 3: expected [88-17=] got [$($c.Range.Text)]"
}
$c = $t.Cell(20,4)
if ($c.Range.Text.StartsWith("63+10=")) {
  $c.Range.Text = "65-7="
} else {
  $failed += "row 20 col 4: expected [63+10=] got [$($c.Range.Text)]"
}
$c = $t.Cell(20,5)
if ($c.Range.Text.StartsWith("3+36=")) {
  $c.Range.Text = "26+66="
} else {
  $failed += "row 20 col 5: expected [3+36=] got [$($c.Range.Text)]"
}

if ($failed.Count -gt 0) {
  Write-Output "FAILED CELLS:"
  foreach ($f in $failed) { Write-Output $f }
} else {
  Write-Output "All cell updates succeeded."
}
